$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 8 (1-indexed, "Numero"=7): Textura cell -> "TileSet" (with spell-check markers)
$cellTextura7 = $t.Cell(8, 2)
$xmlTextura7 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005E0DBE" w:rsidRDefault="00AA7E50" w:rsidP="005E0DBE"><w:pPr><w:jc w:val="center"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>TileSet</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cellTextura7.Range.InsertXML($xmlTextura7)

# Row 8 (1-indexed, "Numero"=7): Comentario cell -> "Prototipo(50x99)"; bookmark _GoBack removed here
$cellComentario7 = $t.Cell(8, 3)
$xmlComentario7 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005E0DBE" w:rsidRDefault="00AA7E50" w:rsidP="005E0DBE"><w:pPr><w:jc w:val="center"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>Prototipo(50x99)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cellComentario7.Range.InsertXML($xmlComentario7)

# Row 9 (1-indexed, "Numero"=8): Textura cell -> "Player Cazador" + bookmark _GoBack (moved here)
$cellTextura8 = $t.Cell(9, 2)
$xmlTextura8 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005E0DBE" w:rsidRDefault="005E0DBE" w:rsidP="005E0DBE"><w:pPr><w:jc w:val="center"/><w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>Player Cazador</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$cellTextura8.Range.InsertXML($xmlTextura8)
